$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.818.63'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '1.637.51'
$ws.Range("E3").Value = '  -1.38%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '308.50'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  -0.11%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3853'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.46%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3807'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.92%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '50.51'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.41%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.321'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.63%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.28%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08351'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.82%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '23.64'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.01%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.933'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.17%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.765'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.78%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.00001304'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").Value = '1.638.69'
$ws.Range("E17").Value = '  -1.24%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '93.31'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06927'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.01%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '19.33'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.04%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.870'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.94%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '13.48'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").Value = '23.825.01'
$ws.Range("E24").Value = '  -0.81%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.430'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.32%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.870'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -8.93%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '21.76'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.71%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '152.90'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.67%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.471'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.18%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '136.29'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.21%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.755'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.95%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.477'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").Value = '1.821.03'
$ws.Range("E33").Value = '  -1.11%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.07941'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.77%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.9786'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -6.36%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.02878'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -4.63%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.551'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.28%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.2647'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.53%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '10.39'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -7.63%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.09054'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.26%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.7466'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.86%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '13.20'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -4.17%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.414'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.31%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '16.53'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.15%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.6850'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.64%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.396'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.46%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.065'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.85%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.04%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.08210'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '133.82'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.14%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.213'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.29%  '
